# Add "remove_data" related mapping rows (debtor.product_group/placement and
# debtor.card_no/card_no) to the config mapping sheet, matching the same
# formatting already used by the existing body rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the two new mapping rows under the existing table (rows 1-16 already
# populated; new data goes into rows 17-18).
$ws.Range("A17").Value = "debtor.product_group"
$ws.Range("B17").Value = "placement"
$ws.Range("A18").Value = "debtor.card_no"
$ws.Range("B18").Value = "card_no"

# Match the style already used for the non-wrapped data cells (columns B for
# the last few existing rows / the header row), so the new cells share the
# same format class instead of the heavier "wrap text" body style.
$ws.Range("A17:B18").WrapText = $false
$ws.Range("A17:B18").Locked = $true

# Move the active selection the way it ended up after the edit.
$ws.Range("C23").Select()

Write-Output "added remove_data mapping rows"
